$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@(76, "InBento", "We do not knowingly collect personally identifiable information from anyone under the age of 18.", 0)
  ,@(77, "InBento", "If you are a parent or guardian and you are aware that your Child has provided us with Personal Data, please contact us", 1)
  ,@(78, "InBento", "If we become aware that we have collected Personal Data from children without verification of parental consent, we take steps to remove that information from our servers", 1)
  ,@(79, "123 Kids Academy", "We do not sell the personal information of Consumers We actually know are less than 16 years of age", 0)
  ,@(80, "123 Kids Academy", "unless We receive affirmative authorization (the 'right to opt-in') from either the Consumer who is between 13 and 16 years of age, or the parent or guardian of a Consumer less than 13 years of age", 0)
  ,@(81, "123 Kids Academy", "If You have reason to believe that a child under the age of 13 (or 16) has provided Us with personal information, please contact Us with sufficient detail to enable Us to delete that information", 1)
  ,@(82, "123 Kids Academy", "Based on the applicable laws of your country, you may have the right to request access to the personal information we collect from you, change that information, or delete it in some circumstances. To request to review, update, or delete your personal information", 1)
  ,@(83, "iCivics", "Consistent with the Children's Online Privacy Protection Act of 1998 ('COPPA'), we will never knowingly request personally identifying information from anyone under the age of 13 without prior verifiable parental consent", 0)
  ,@(84, "iCivics", "If a child under the age of 13 wishes to sign up for a user account on our Website, he or she must first obtain parental approval and provide a parent's email address for confirmation", 0)
  ,@(85, "iCivics", "If you provide iCivics with personal information, you may access and edit this information at any time by visiting the “Edit Account” link at the top right corner of the Website", 1)
  ,@(86, "iCivics", "If you would like to request deletion of your account, please email us at feedback@iCivics.org", 1)
  ,@(87, "Chuggington", "We may collect information such as your first and last name, home or other physical address, email address, telephone number, date of birth and other information that permits you to be contacted physically or online, and other information, including information about your interactions with us or others", 0)
  ,@(88, "Chuggington", "You may also provide us information about others, e.g., when purchasing tickets or registering for an event", 0)
  ,@(89, "Chuggington", "We may collect payment information when you transact with us", 0)
  ,@(90, "Chuggington", "you have the right to receive: a) information identifying any third party company(ies) to whom we may have disclosed, within the previous calendar year", 1)
  ,@(91, "Chuggington", "your Personally Identifiable Information for that company's direct marketing purposes", 1)
  ,@(92, "Chuggington", "a description of the categories of Personally Identifiable Information disclosed", 1)
  ,@(93, "Chuggington", "To obtain such information, please email your request to info@hfecorp.com", 1)
  ,@(94, "Sesame Workshop", "Visitors can provide their first name, a friend’s first name and friend’s email address to send a feature of a Children’s Platform to the friend via email", 0)
  ,@(95, "Sesame Workshop", "We use such information only to fulfill the visitor’s request and do not store that information for any longer than reasonably necessary to send the relevant email", 0)
  ,@(96, "Sesame Workshop", "When visitors on our Children’s Platforms wish to sign up to receive our email newsletters, we will ask for a parent’s first name and email address", 0)
  ,@(97, "Sesame Workshop", "We will send a notification email to the parent and give the parent the option to decline the newsletter", 0)
  ,@(98, "Sesame Workshop", "We do not make personal information collected from children publicly available nor do we enable children to so on our Platforms", 0)
  ,@(99, "Sesame Workshop", "Parents may contact us to review, update or delete any of their children’s personal information that we may have collected and to elect for us not to collect any additional personal information from their children", 1)
  ,@(100, "Sesame Workshop", "Parents may contact us (at dotorgmail@sesame.org)", 1)
  ,@(101, "Amaya Kids", "Child users do not have the ability to provide any personal information when using the Games beyond what has collected automatically (as noted above)", 0)
  ,@(102, "Amaya Kids", "We do not make personal information collected from children publicly available nor do we enable children to so on our Platforms", 0)
  ,@(103, "Mattel", "We do not collect personal contact information from children at Mattel Services directed to children without the consent of a parent or legal guardian, except in limited circumstances authorized by law", 0)
  ,@(104, "Mattel", "We do not knowingly sell the personal information of minors under 16 years of age", 0)
  ,@(105, "Sagomini", "We collect non-personal data indirectly from the following third-party analytics providers", 0)
  ,@(106, "Sagomini", "These third parties do not share your user data; they are used to support the operation of our apps", 0)
  ,@(107, "Sagomini", "If you want to opt out of services, or review or delete you or your child’s information, please contact us at privacy@sagosago.com", 1)
  ,@(108, "Scratch Jr", "Where applicable, we indicate whether and why you must provide us with your Personal Information, as well as the consequences of failing to do so", 0)
  ,@(109, "Scratch Jr", "Deleting your information. You can delete your projects by deleting the App from your device", 1)
  ,@(110, "Crayola", "To play our apps we do not collect any personal data from you", 0)
  ,@(111, "Crayola", "Right of Deletion - You can ask a company to delete any information they may have about you", 1)
  ,@(112, "BabyBus", "Email Address To create your account for the App as per your request With your consent", 0)
  ,@(113, "BabyBus", "You have the right to delete your personal information", 1)
  ,@(114, "Path of Giants", "The Service is not directed to children under the age of 13", 0)
  ,@(115, "Path of Giants", "we allow children below the age of 13 to use the Service, we do not knowingly collect personal information from children under the age of 13 without first obtaining verifiable parental consent", 0)
  ,@(116, "Path of Giants", "If we learn that we inadvertently collected personal information from a children under the age of 13 without first obtaining verifiable parental consent, we will delete that information as quickly as possible", 1)
  ,@(117, "Path of Giants", "If you are a parent or guardian of a child who you believe provided Journey Bound with personal information without your consent, please contact us at contact@journeyboundgames.com", 1)
  ,@(118, "Wollstonecraft", "We do not knowingly collect personal information from children under the age of 13 yeras old", 0)
  ,@(119, "Wollstonecraft", "If a user is identified as under 13 years old, we will not collect or use any information of this user and we will delete any information already recieved in a secure manner", 1)
  ,@(120, "Wollstonecraft", "If you have any questions about our privacy policy, contact us by email at info@hololabs.org", 1)
  ,@(121, "Prodigy", "We do not disclose the personal information of Students to third parties for marketing or promotional purposes", 0)
  ,@(122, "Prodigy", "We use all User information we collect for the following purposes: To monitor, maintain, analyze and improve functionality of Services", 0)
  ,@(123, "Prodigy", "may contact our customer support at any time to request access to, deletion of, or correction of any personal information we have collected from or about Student Users associated with their account", 1)
  ,@(124, "Prodigy", "including a request to us to cease collecting personal information from those Student Users", 1)
  ,@(125, "Prodigy", "For information on how consumers can ask questions or file complaints related to Prodigy's Privacy Policy and practices, please email COPPAPrivacy@ikeepsafe.org", 1)
  ,@(126, "Balloon Crush", "We do not knowingly collect personally identifiable information from anyone under the age of 13", 0)
  ,@(127, "Balloon Crush", "Our Service may contain links to other websites that are not operated by Us", 0)
  ,@(128, "Balloon Crush", "If You click on a third party link, You will be directed to that third party's site. We strongly advise You to review the Privacy Policy of every site You visit", 0)
  ,@(129, "Balloon Crush", "If We become aware that We have collected Personal Data from anyone under the age of 13 without verification of parental consent, We take steps to remove that information from Our servers", 1)
  ,@(130, "Balloon Crush", "If you have any questions about this privacy Policy, You can contact us: By email: pradyuishere@gmail.com", 1)
  ,@(131, "Crescent Moon Games", "We do not knowingly collect personally identifiable information from children under 13", 0)
  ,@(132, "Crescent Moon Games", "In the case we discover that a child under 13 has provided us with personal information, we immediately delete this from our servers", 1)
  ,@(133, "Crescent Moon Games", "If you are a parent or guardian and you are aware that your child has provided us with personal information, please contact us so that we will be able to do necessary actions", 1)
  ,@(134, "Adventuring Academy", "from Child Users, we may collect information about patterns of usage such as which activities a child commences and completes, when a child starts and stops an activity, and which areas of the Services the child frequents", 0)
  ,@(135, "Adventuring Academy", "The information collected through these technical methods on the child-directed portions of the Services are used only to support the internal operations of the Services", 0)
  ,@(136, "Adventuring Academy", "In addition, Adult Users may contact us at any time as described in Section 13 (Contact Us) below to request that we provide for their review, or delete from our records, any PI they have provided about Child Users associated with their Accounts, or to cease collecting PI from those Child Users, as applicable", 1)
  ,@(137, "Animal Jam", "We collect what is reasonably necessary for us to provide children with access to Animal Jam, such as a username, password, demographic information, and a parental email address", 0)
  ,@(138, "Animal Jam", "For example, if you or your child submits Information in order to create an account with us, or otherwise contacts us directly", 0)
  ,@(139, "Animal Jam", "You may also request that we no longer collect Personal Information from your child, or have us delete your child's Personal Information we have collected", 1)
  ,@(140, "IDZ Digital", "We do not knowingly collect or solicit personal data about or direct or target interest-based advertising to anyone under the age of 13 or knowingly allow such persons to use our Services", 0)
  ,@(141, "IDZ Digital", "If we learn that we have collected personal data about a child under age 13, we will delete that data as quickly as possible", 1)
  ,@(142, "IDZ Digital", "If a parent/guardian becomes aware that a child under the age of 13 has attempted to contact and/or do business with us, please advise us by email at support@idzdigital.com so that we may rectify the situation", 1)
  ,@(143, "Infinity Games", "We do not knowingly collect personally identifiable information from children under 16", 0)
  ,@(144, "Infinity Games", "If we become aware that we have collected Personal Information from a child under age 16 without verification of parental consent, we will take steps to remove that information from our servers", 1)
  ,@(145, "Infinity Games", "If you are a parent or guardian and you are aware that your Children has provided us with Personal Information, please contact us", 1)
  ,@(146, "E One", "Accordingly, we do not generally expect to collect personal information directly from any children, except at the direction of a parent/guardian", 0)
  ,@(147, "E One", "Limited children’s information if you permit, including name(s), date(s) of birth and gender, as well as any information you (or your child) submit as part of a promotion or competition", 0)
  ,@(148, "E One", "If you believe we have collected personal information from a child under 13 without parental consent, please let us know via the contact information below and we will endeavor to promptly delete it", 1)
  ,@(149, "E One", "If you have such concerns, we request that you initially contact us (using the contact details below) so that we can investigate, and hopefully resolve, your concerns", 1)
  ,@(150, "SpeedyMind", "We follow the COPPA (Children’s Online Privacy Protection Act) compliant privacy practices. We do not collect any identifying data from kids", 0)
  ,@(151, "SpeedyMind", "You have the right to request the erase your Personal Information under certain conditions of this Policy", 1)
  ,@(152, "SpeedyMind", "If you wish to be informed what Personal Information we hold about you and if you want it to be removed from our systems, please contact us", 1)
  ,@(153, "Budge Studios", "Budge Studios does not knowingly collect personal information as defined by COPPA from children through its Sites or Programs", 0)
  ,@(154, "Budge Studios", "We collect, use and share personal information from children through our Apps as follows:", 0)
  ,@(155, "Budge Studios", "Additionally, any parent has the right to: (1) review, correct, or delete the child's personal information; and/or (2) discontinue further collection, use, or sharing of the child's personal information", 1)
  ,@(156, "Budge Studios", "If you are a parent or guardian and believe that we have collected your child's personal information in violation of COPPA, please contact us as set out in the How Do I Contact Budge Studios? Section, and we will remove such information to the extent required by COPPA", 1)
  ,@(157, "Amanita Design", "We do not collect, process or store any personal data from any of our games on any platform", 0)
  ,@(158, "Amanita Design", "Parents can learn more about how we handle data by contacting us at email address below", 1)
  ,@(159, "NjoyKidz", "We do not ask for or collect any child-specific data like concrete name and surname, location, school, etc., anywhere on the Platform", 0)
  ,@(160, "NjoyKidz", "We may ask for your child’s interests but that is just to personalize your child’s experience on the Platform and show them the relevant and appropriate content", 0)
  ,@(161, "Yateland", "We don’t knowingly collect personal information from children under the age of 13 in violation of COPPA, and if in the event that a user identifies himself or herself as a child under the age of 13 through a support request", 0)
  ,@(162, "Yateland", "we will not collect, store or use, and will delete in a secure manner, any personal information of such user", 0)
  ,@(163, "American Heart Association", "With parental consent, we may collect information from children under the age of 13 such as: name, address, email address, account information, school, messages sent to us through our chat interface, and content they create themselves", 0)
  ,@(164, "American Heart Association", "Children under the age of 13 may be able to make certain content such as content they create themselves visible to others or the public", 0)
  ,@(165, "American Heart Association", "This might include, for example, a webpage or parts of webpages operated by children that have been designed for group or public viewing, or photos of themselves involved in our programs or other activities", 0)
  ,@(166, "American Heart Association", "Regardless of what is displayed or submitted, parents can revoke their consent, request that information about their children be hidden or, in some cases, deleted", 1)
  ,@(167, "American Heart Association", "Request that information about their children be hidden or, in some cases, deleted, by contacting our offices by phone at the number at the bottom of this page or via email", 1)
  ,@(168, "IT SYSTEM", "We do not collect or require users to enter their personal information when using our Products. We do not collect any personal information from children with our Products", 0)
  ,@(169, "IT SYSTEM", "When a user identifies himself or herself as a child under the age of 13 through a support request or through any feedback, we will not collect, store or use, and will delete in a secure manner, any personal information of such user", 1)
  ,@(170, "Pazu Games", "We do not knowingly collect personally identifiable information from anyone under the age of 18", 0)
  ,@(171, "Pazu Games", "If we become aware that we have collected Personal Data from children without verification of parental consent, we take steps to remove that information from our servers", 1)
  ,@(172, "Pazu Games", "If you are a parent or guardian and you are aware that your Child has provided us with Personal Data, please contact us", 1)
  ,@(173, "Oleg", "I do not knowingly collect personally identifiable information from children under 13 years of age", 0)
  ,@(174, "Oleg", "If you are a parent or guardian and you are aware that your child has provided us with personal information, please contact me so that I will be able to do the necessary actions", 1)
  ,@(175, "Grapefrukt", "No user-identifiable data is kept", 0)
)

# Pre-seed shared-string insertion order for two pairs whose original
# insertion order (in the source workbook) differs from row-scan order:
# the text below is registered in the shared-string table right after
# column A of the row is written, but before column B of that same row,
# matching the table order recorded in the target workbook.
$preseedAfterA = @{
  88 = "We may collect payment information when you transact with us"
  126 = "Our Service may contain links to other websites that are not operated by Us"
}

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  if ($preseedAfterA.ContainsKey($r)) {
    $ws.Range("Z1").Value = $preseedAfterA[$r]
    $ws.Range("Z1").ClearContents()
  }
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
}

$ws.Range("C175").Select() | Out-Null